$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.727.51'
$ws.Range("E2").Value = '  +0.40%  '
# Row 3
$ws.Range("D3").Value = '3.523.16'
$ws.Range("E3").Value = '  +0.91%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.88'
$ws.Range("E5").Value = '  -0.35%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.14'
$ws.Range("E6").Value = '  +5.45%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  +0.35%  '
# Row 8
$ws.Range("E8").Value = '  -0.05%  '
# Row 9
$ws.Range("E9").Value = '  -6.74%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.649'
$ws.Range("E10").Value = '  -0.18%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.76'
$ws.Range("E11").Value = '  +1.19%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000302'
$ws.Range("E12").Value = '  -2.26%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.51'
$ws.Range("E13").Value = '  -0.20%  '
# Row 14
$ws.Range("D14").Value = '4.083.27'
$ws.Range("E14").Value = '  +1.23%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '595.89'
$ws.Range("E15").Value = '  -1.36%  '
# Row 16
$ws.Range("E16").Value = '  +1.45%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.16'
$ws.Range("E17").Value = '  +1.47%  '
# Row 18
$ws.Range("D18").Value = '69.923.52'
$ws.Range("E18").Value = '  +0.63%  '
# Row 19
$ws.Range("D19").Value = '3.526.03'
$ws.Range("E19").Value = '  -2.10%  '
# Row 20
$ws.Range("E20").Value = '  +1.40%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.992'
$ws.Range("E21").Value = '  +0.45%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.11'
$ws.Range("E22").Value = '  +5.68%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.31'
$ws.Range("E23").Value = '  +5.16%  '
# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.35'
$ws.Range("E24").Value = '  -3.37%  '
# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.66'
$ws.Range("E25").Value = '  +0.45%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.17'
$ws.Range("E26").Value = '  +4.74%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.87'
$ws.Range("E27").Value = '  -0.64%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.61'
$ws.Range("E28").Value = '  -1.35%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.51'
$ws.Range("E29").Value = '  -0.39%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  +1.11%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.28'
$ws.Range("E31").Value = '  +5.41%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.46'
$ws.Range("E32").Value = '  +0.13%  '
# Row 33
$ws.Range("E33").Value = '  +0.04%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.14'
$ws.Range("E34").Value = '  -0.17%  '
# Row 35
$ws.Range("D35").Value = '0.0₃0852'
$ws.Range("E35").Value = '  +9.53%  '
# Row 36
$ws.Range("D36").Value = '3.719.54'
$ws.Range("E36").Value = '  +3.00%  '
# Row 37
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.02%  '
# Row 38
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.07'
$ws.Range("E38").Value = '  -3.74%  '
# Row 39
$ws.Range("E39").Value = '  +0.18%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.393'
$ws.Range("E40").Value = '  -1.20%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.57'
$ws.Range("E41").Value = '  -0.47%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '490.44'
$ws.Range("E42").Value = '  -6.58%  '
# Row 43
$ws.Range("E43").Value = '  -3.58%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0455'
$ws.Range("E44").Value = '  -0.77%  '
# Row 45
$ws.Range("E45").Value = '  -1.49%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("E46").Value = '  -4.45%  '
# Row 47
$ws.Range("E47").Value = '  -1.46%  '
# Row 48
$ws.Range("E48").Value = '  +0.34%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.51'
$ws.Range("E49").Value = '  -3.55%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000246'
$ws.Range("E50").Value = '  +1.63%  '
# Row 51
$ws.Range("E51").Value = '  +10.92%  '
